$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Score corrections (row 7 and row 11) ---
$ws.Range("D7").Value = 0
$ws.Range("B11").Value = 3.5
$ws.Range("D11").Value = 1.75

# --- New "Nhận xét" (Remarks) column F ---
$ws.Range("F1").Value = "Nhận xét"
$ws.Range("F1").Borders.LineStyle = 1

$ws.Range("F2").Value = "Chưa đưa ra được thuật toán bài 3"
$ws.Range("F2").Borders.LineStyle = 1

$ws.Range("F3").Value = "OK"
$ws.Range("F3").Borders.LineStyle = 1

$ws.Range("F4").Value = "OK"
$ws.Range("F4").Borders.LineStyle = 1

$ws.Range("F5").Value = "OK"
$ws.Range("F5").Borders.LineStyle = 1

$ws.Range("F6").Value = "OK"
$ws.Range("F6").Borders.LineStyle = 1

$ws.Range("F7").Value = "Chưa đưa ra được thuật toán câu c bài 3, câu b, c bài 1 sai trong quá trình làm dẫn đến kết quả sai"
$ws.Range("F7").Borders.LineStyle = 1

$ws.Range("F8").Value = "OK"
$ws.Range("F8").Borders.LineStyle = 1

$ws.Range("F9").Value = "OK"
$ws.Range("F9").Borders.LineStyle = 1

$ws.Range("F10").Value = "OK"
$ws.Range("F10").Borders.LineStyle = 1

$ws.Range("F11").Value = "Chưa đưa ra được thuật toán câu c bài 3, câu b bài 3 chưa giải thích rõ ràng."
$ws.Range("F11").Borders.LineStyle = 1

$ws.Range("F12").Value = "Chưa đưa ra được thuật toán bài 3"
$ws.Range("F12").Borders.LineStyle = 1

# F13 already carries the stray date-number-format style (s=2 in the
# original file); just add the border to it, preserving that format, then
# set its text.
$ws.Range("F13").Borders.LineStyle = 1
$ws.Range("F13").Value = "OK"

$ws.Range("F14").Value = "OK"
$ws.Range("F14").Borders.LineStyle = 1

$ws.Range("F15").Value = "Chưa làm được câu 3"
$ws.Range("F15").Borders.LineStyle = 1

$ws.Range("F16").Value = "OK"
$ws.Range("F16").Borders.LineStyle = 1

$ws.Range("F17").Value = "OK"
$ws.Range("F17").Borders.LineStyle = 1

# --- Extend the merged "Link bài làm các nhóm" footer row to column F ---
$ws.Range("A18:F18").Merge()
# Re-apply the uniform thin box border/centered style across the whole
# merged region (merging alone redistributes the old A18:E18 border into
# mismatched per-cell edges; this restores one consistent look matching
# the rest of the row, now also covering the new F18 cell).
$ws.Range("A18:F18").Borders.LineStyle = 1

# --- Selection moves to F10 ---
$ws.Range("F10").Select()
